$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.846.72"
$ws.Range("E2").Value = "  -3.29%  "
$ws.Range("D3").Value = "2.254.66"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.54%  "
$ws.Range("D5").Value = "'493.91"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").Value = "2.285.72"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").Value = "'0.0937"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'0.321"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'4.61"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("D14").Value = "2.666.06"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "'21.50"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "53.795.97"
$ws.Range("E16").Value = "  -3.30%  "
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "2.318.34"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").Value = "'9.88"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'302.21"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("D22").Value = "'6.41"
$ws.Range("E22").Value = "  +3.04%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'5.34"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("D25").Value = "'63.70"
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'0.373"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").Value = "2.394.11"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").Value = "'0.148"
$ws.Range("E29").Value = "  +2.05%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").Value = "'165.28"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").Value = "'1.59"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").Value = "'17.57"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("D40").Value = "'0.867"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("D41").Value = "'3.61"
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").Value = "'35.40"
$ws.Range("D43").Value = "'0.374"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("D45").Value = "'3.34"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").Value = "'125.76"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").Value = "'0.544"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'236.12"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  +1.35%  "
